$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 803, pushing the existing rows 803:844 down to 804:845
$ws.Rows("803:803").Insert()

# Populate the newly inserted row with the new data point.
# The leading apostrophe forces the date-looking text to be stored as
# literal text (matching how the rest of column A is stored), instead of
# being auto-converted into a real Excel date serial value.
$ws.Range("A803").Value = "'2026/02/13"
$ws.Range("B803").Value = "金"
$ws.Range("C803").Value = 10
$ws.Range("D803").Value = 201
